$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some of the new Price values (column D) look like plain numbers to Excel
# (e.g. "314.29"), which would otherwise be silently auto-converted into a
# numeric value instead of being kept as text, losing the original
# formatting (trailing zeros, exact digit count, etc). Force those specific
# cells to Text format before assigning the value, then restore the default
# "Normal" style afterwards so no stray number format is left behind.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

$ws.Range("D5").Value = '314.29'
$ws.Range("D6").Value = '1.001'
$ws.Range("D7").Value = '0.4473'
$ws.Range("D9").Value = '0.07517'
$ws.Range("D10").Value = '0.8929'
$ws.Range("D13").Value = '6.758'
$ws.Range("D14").Value = '94.33'
$ws.Range("D15").Value = '5.414'
$ws.Range("D16").Value = '0.07115'
$ws.Range("D17").Value = '1.002'
$ws.Range("D18").Value = '0.000008813'
$ws.Range("D23").Value = '10.93'
$ws.Range("D25").Value = '1.974'
$ws.Range("D26").Value = '2.380'
$ws.Range("D27").Value = '151.43'
$ws.Range("D29").Value = '5.372'
$ws.Range("D30").Value = '117.64'
$ws.Range("D31").Value = '0.08832'
$ws.Range("D32").Value = '0.7856'
$ws.Range("D33").Value = '1.204'
$ws.Range("D34").Value = '4.523'
$ws.Range("D35").Value = '2.891'
$ws.Range("D36").Value = '1.000'
$ws.Range("D38").Value = '0.01992'
$ws.Range("D39").Value = '0.05333'
$ws.Range("D40").Value = '7.395'
$ws.Range("D41").Value = '0.5325'
$ws.Range("D42").Value = '0.1734'
$ws.Range("D43").Value = '2.861'
$ws.Range("D44").Value = '2.289'
$ws.Range("D45").Value = '8.784'
$ws.Range("D46").Value = '0.5128'
$ws.Range("D47").Value = '10.77'
$ws.Range("D48").Value = '106.12'
$ws.Range("D49").Value = '1.704'

$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"

# Remaining cells (percentage strings, multi-dot price strings, etc.) are
# never mistaken for numbers by Excel, so a plain .Value assignment keeps
# them as text already.
$ws.Range("D2").Value = '27.405.93'
$ws.Range("E2").Value = '  +1.35%  '
$ws.Range("D3").Value = '1.827.25'
$ws.Range("E3").Value = '  -0.08%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("E5").Value = '  +0.84%  '
$ws.Range("E6").Value = '  +0.03%  '
$ws.Range("E7").Value = '  +3.26%  '
$ws.Range("E8").Value = '  +2.50%  '
$ws.Range("E9").Value = '  +3.14%  '
$ws.Range("E10").Value = '  +5.69%  '
$ws.Range("E11").Value = '  +1.87%  '
$ws.Range("D12").Value = '1.826.44'
$ws.Range("E12").Value = '  -0.02%  '
$ws.Range("E13").Value = '  +1.36%  '
$ws.Range("E14").Value = '  +5.10%  '
$ws.Range("E15").Value = '  +2.20%  '
$ws.Range("E16").Value = '  +0.45%  '
$ws.Range("E17").Value = '  +0.00%  '
$ws.Range("E18").Value = '  +0.36%  '
$ws.Range("E20").Value = '  +2.05%  '
$ws.Range("D21").Value = '27.424.08'
$ws.Range("E21").Value = '  +1.30%  '
$ws.Range("E22").Value = '  +2.47%  '
$ws.Range("E23").Value = '  +0.19%  '
$ws.Range("D24").Value = '2.055.93'
$ws.Range("E25").Value = '  -0.36%  '
$ws.Range("E27").Value = '  +0.19%  '
$ws.Range("E28").Value = '  +1.44%  '
$ws.Range("E29").Value = '  +2.59%  '
$ws.Range("E30").Value = '  +0.48%  '
$ws.Range("E31").Value = '  +1.00%  '
$ws.Range("E32").Value = '  +6.13%  '
$ws.Range("E33").Value = '  +2.02%  '
$ws.Range("E34").Value = '  +2.00%  '
$ws.Range("E35").Value = '  -0.29%  '
$ws.Range("E36").Value = '  +0.05%  '
$ws.Range("E37").Value = '  +1.74%  '
$ws.Range("E38").Value = '  +2.42%  '
$ws.Range("E40").Value = '  +2.68%  '
$ws.Range("E41").Value = '  +3.54%  '
$ws.Range("E42").Value = '  +1.87%  '
$ws.Range("E43").Value = '  -0.17%  '
$ws.Range("E44").Value = '  +17.50%  '
$ws.Range("E45").Value = '  +2.52%  '
$ws.Range("E46").Value = '  +7.82%  '
$ws.Range("E47").Value = '  +2.03%  '
$ws.Range("E48").Value = '  +0.33%  '
$ws.Range("E49").Value = '  +2.08%  '
$ws.Range("E50").Value = '  +0.09%  '
$ws.Range("E51").Value = '  +0.62%  '
